$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper: set a cell to a numeric-looking text value while preserving it
# as a text string (matching the workbook's existing inlineStr/text cells)
# instead of letting Excel auto-convert it into a floating point number.
function Set-TextValue($row, $col, $text) {
    $cell = $ws.Cells.Item($row, $col)
    $cell.Value = "'" + $text
    $cell.Style = "Normal"
}

# Column D (Price) updates
Set-TextValue 2  4 "261.35"
Set-TextValue 3  4 "22.92"
Set-TextValue 4  4 "6.203"
Set-TextValue 6  4 "6.731"
Set-TextValue 7  4 "3.458"
Set-TextValue 8  4 "1.344"
Set-TextValue 9  4 "0.7990"
Set-TextValue 10 4 "0.1588"
Set-TextValue 11 4 "0.08109"
Set-TextValue 14 4 "0.09309"
Set-TextValue 15 4 "3.863"
Set-TextValue 16 4 "0.001700"
Set-TextValue 17 4 "0.04794"
Set-TextValue 18 4 "0.01335"
Set-TextValue 20 4 "0.001093"
Set-TextValue 22 4 "0.0001500"
Set-TextValue 24 4 "2.214"
Set-TextValue 27 4 "0.0003202"
Set-TextValue 40 4 "0.04611"
Set-TextValue 41 4 "0.007180"
Set-TextValue 42 4 "0.1117"
Set-TextValue 43 4 "0.003600"
Set-TextValue 44 4 "0.01022"
Set-TextValue 45 4 "0.002972"
Set-TextValue 46 4 "0.00005920"
Set-TextValue 48 4 "0.7001"
Set-TextValue 49 4 "0.1561"
Set-TextValue 50 4 "0.00002100"

# Column E (Volume(1h)) text updates
$ws.Cells.Item(18, 5).Value = "17OneONE"
$ws.Cells.Item(21, 5).Value = "20HotbitTokenHTBWorstin24h"
